$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HOME")

# -- New "Home section 6" rows (61-65) --------------------------------

# Row 61
$ws.Range("A61").Value = "SECTION_6_TEXT_1"
$ws.Range("B61").Value = "hmmm" + [char]0x2026

# Row 62 (label first, value filled in later to match authoring order)
$ws.Range("A62").Value = "SECTION_6_TEXT_2"

# Row 63
$ws.Range("A63").Value = "SECTION_6_TEXT_3"
$ws.Range("B63").Value = "Hello, Neko.`nTell me what you can do!"
$ws.Range("B63").WrapText = $true
$ws.Rows.Item(63).RowHeight = 30

# Row 64 (label first, value filled in later to match authoring order)
$ws.Range("A64").Value = "SECTION_6_TEXT_4"

# Row 65
$ws.Range("A65").Value = "SECTION_6_TEXT_5"
$ws.Range("B65").Value = "Hey Neko. Surprised me!"

# Back-fill B62 / B64 (this is the order the shared strings were
# originally authored in)
$ws.Range("B62").Value = "let's`nsee"
$ws.Range("B62").WrapText = $true
$ws.Rows.Item(62).RowHeight = 30

$ws.Range("B64").Value = "Hi, Neko.`nI have many ideas. But I don't know`nwhere to start..."
$ws.Range("B64").WrapText = $true
$ws.Rows.Item(64).RowHeight = 45

# -- View state: scroll down and select B65 ---------------------------
[void]$ws.Range("B65").Select()
$excel.ActiveWindow.ScrollRow = 55
